# Update the simulated log-write-mode results for data set gr25_06.
# Refreshes run_time, max_er, and the per-iteration convergence metrics
# (iter 0..iter 19) for each of the 10 simulated runs (rows 2-11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.2279942035675049
$ws.Cells.Item(2, 5).Value = 56.10687993525971
$ws.Cells.Item(2, 6).Value = 0.002183817042338903
$ws.Cells.Item(2, 7).Value = 0.001779179316321811
$ws.Cells.Item(2, 8).Value = 0.001623303116226225
$ws.Cells.Item(2, 9).Value = 0.001552078112251966
$ws.Cells.Item(2, 10).Value = 0.001475118993924459
$ws.Cells.Item(2, 11).Value = 0.001448279285121114
$ws.Cells.Item(2, 12).Value = 0.0014041121637863
$ws.Cells.Item(2, 13).Value = 0.001326784426068892
$ws.Cells.Item(2, 14).Value = 0.001326784426068892
$ws.Cells.Item(2, 15).Value = 0.001293642194379384
$ws.Cells.Item(2, 16).Value = 0.001243134721200023
$ws.Cells.Item(2, 17).Value = 0.001240165102122812
$ws.Cells.Item(2, 18).Value = 0.001240165102122812
$ws.Cells.Item(2, 19).Value = 0.001204690709905671
$ws.Cells.Item(2, 20).Value = 0.001154922463341095
$ws.Cells.Item(2, 21).Value = 0.001154922463341095
$ws.Cells.Item(2, 22).Value = 0.001136975599171764
$ws.Cells.Item(2, 23).Value = 0.001105598984748099
$ws.Cells.Item(2, 24).Value = 0.00109787980583268
$ws.Cells.Item(2, 25).Value = 0.001093701363260423

$ws.Cells.Item(3, 3).Value = 0.2160384654998779
$ws.Cells.Item(3, 5).Value = 53.96451203334982
$ws.Cells.Item(3, 6).Value = 0.002183817042338903
$ws.Cells.Item(3, 7).Value = 0.001747333769651048
$ws.Cells.Item(3, 8).Value = 0.00154987233688212
$ws.Cells.Item(3, 9).Value = 0.001456229798269619
$ws.Cells.Item(3, 10).Value = 0.001394005360621007
$ws.Cells.Item(3, 11).Value = 0.001368221418100373
$ws.Cells.Item(3, 12).Value = 0.001292530709161425
$ws.Cells.Item(3, 13).Value = 0.001139396068373133
$ws.Cells.Item(3, 14).Value = 0.001139396068373133
$ws.Cells.Item(3, 15).Value = 0.001139396068373133
$ws.Cells.Item(3, 16).Value = 0.001139396068373133
$ws.Cells.Item(3, 17).Value = 0.001139396068373133
$ws.Cells.Item(3, 18).Value = 0.001095653656687426
$ws.Cells.Item(3, 19).Value = 0.001095653656687426
$ws.Cells.Item(3, 20).Value = 0.001095653656687426
$ws.Cells.Item(3, 21).Value = 0.001079275579990986
$ws.Cells.Item(3, 22).Value = 0.001068592176468574
$ws.Cells.Item(3, 23).Value = 0.001058460986854091
$ws.Cells.Item(3, 24).Value = 0.001058460986854091
$ws.Cells.Item(3, 25).Value = 0.00105193980571832

$ws.Cells.Item(4, 3).Value = 0.2080078125
$ws.Cells.Item(4, 5).Value = 54.77745932613107
$ws.Cells.Item(4, 6).Value = 0.002183817042338903
$ws.Cells.Item(4, 7).Value = 0.001770272123480519
$ws.Cells.Item(4, 8).Value = 0.001537293367222318
$ws.Cells.Item(4, 9).Value = 0.001513367714587996
$ws.Cells.Item(4, 10).Value = 0.001513367714587996
$ws.Cells.Item(4, 11).Value = 0.001513367714587996
$ws.Cells.Item(4, 12).Value = 0.001424643295353621
$ws.Cells.Item(4, 13).Value = 0.001196413938546844
$ws.Cells.Item(4, 14).Value = 0.001196413938546844
$ws.Cells.Item(4, 15).Value = 0.001159350172556639
$ws.Cells.Item(4, 16).Value = 0.001159350172556639
$ws.Cells.Item(4, 17).Value = 0.001159350172556639
$ws.Cells.Item(4, 18).Value = 0.001159350172556639
$ws.Cells.Item(4, 19).Value = 0.001148774925146554
$ws.Cells.Item(4, 20).Value = 0.00113328906723255
$ws.Cells.Item(4, 21).Value = 0.001125900715091843
$ws.Cells.Item(4, 22).Value = 0.001090747149171779
$ws.Cells.Item(4, 23).Value = 0.001079538471415672
$ws.Cells.Item(4, 24).Value = 0.001079538471415672
$ws.Cells.Item(4, 25).Value = 0.001067786731503529

$ws.Cells.Item(5, 3).Value = 0.2785139083862305
$ws.Cells.Item(5, 5).Value = 54.67943350428141
$ws.Cells.Item(5, 6).Value = 0.002162077097587758
$ws.Cells.Item(5, 7).Value = 0.001737696203857617
$ws.Cells.Item(5, 8).Value = 0.001737696203857617
$ws.Cells.Item(5, 9).Value = 0.001589896623294049
$ws.Cells.Item(5, 10).Value = 0.001565168247107918
$ws.Cells.Item(5, 11).Value = 0.001438931452340169
$ws.Cells.Item(5, 12).Value = 0.001331488087718228
$ws.Cells.Item(5, 13).Value = 0.001331488087718228
$ws.Cells.Item(5, 14).Value = 0.00125998331238809
$ws.Cells.Item(5, 15).Value = 0.001196629534628492
$ws.Cells.Item(5, 16).Value = 0.001168215956809665
$ws.Cells.Item(5, 17).Value = 0.001092997837192642
$ws.Cells.Item(5, 18).Value = 0.001092997837192642
$ws.Cells.Item(5, 19).Value = 0.001092997837192642
$ws.Cells.Item(5, 20).Value = 0.001092997837192642
$ws.Cells.Item(5, 21).Value = 0.001092997837192642
$ws.Cells.Item(5, 22).Value = 0.001092997837192642
$ws.Cells.Item(5, 23).Value = 0.001076957324023832
$ws.Cells.Item(5, 24).Value = 0.001066431138180037
$ws.Cells.Item(5, 25).Value = 0.001065875896769618

$ws.Cells.Item(6, 3).Value = 0.2839963436126709
$ws.Cells.Item(6, 5).Value = 60.26204843626874
$ws.Cells.Item(6, 6).Value = 0.002183817042338903
$ws.Cells.Item(6, 7).Value = 0.001736917710231314
$ws.Cells.Item(6, 8).Value = 0.001700336605672082
$ws.Cells.Item(6, 9).Value = 0.001560510874294788
$ws.Cells.Item(6, 10).Value = 0.001480708646370181
$ws.Cells.Item(6, 11).Value = 0.001480708646370181
$ws.Cells.Item(6, 12).Value = 0.001480708646370181
$ws.Cells.Item(6, 13).Value = 0.001431665875917394
$ws.Cells.Item(6, 14).Value = 0.001400990202078727
$ws.Cells.Item(6, 15).Value = 0.001364412942884014
$ws.Cells.Item(6, 16).Value = 0.001362624379046705
$ws.Cells.Item(6, 17).Value = 0.001293812999451367
$ws.Cells.Item(6, 18).Value = 0.001248775592042907
$ws.Cells.Item(6, 19).Value = 0.001248775592042907
$ws.Cells.Item(6, 20).Value = 0.001220148915292217
$ws.Cells.Item(6, 21).Value = 0.001220148915292217
$ws.Cells.Item(6, 22).Value = 0.001181760048829812
$ws.Cells.Item(6, 23).Value = 0.001181760048829812
$ws.Cells.Item(6, 24).Value = 0.001180135588875932
$ws.Cells.Item(6, 25).Value = 0.001174698799927266

$ws.Cells.Item(7, 3).Value = 0.2420375347137451
$ws.Cells.Item(7, 5).Value = 64.18114163232713
$ws.Cells.Item(7, 6).Value = 0.002164975666093337
$ws.Cells.Item(7, 7).Value = 0.001880346008562489
$ws.Cells.Item(7, 8).Value = 0.001629091367249202
$ws.Cells.Item(7, 9).Value = 0.001629091367249202
$ws.Cells.Item(7, 10).Value = 0.001588984512539374
$ws.Cells.Item(7, 11).Value = 0.001541930049308511
$ws.Cells.Item(7, 12).Value = 0.001465184356307575
$ws.Cells.Item(7, 13).Value = 0.001404741955364861
$ws.Cells.Item(7, 14).Value = 0.001389983832468538
$ws.Cells.Item(7, 15).Value = 0.001346302062314716
$ws.Cells.Item(7, 16).Value = 0.001346302062314716
$ws.Cells.Item(7, 17).Value = 0.001322526218484331
$ws.Cells.Item(7, 18).Value = 0.001292322669743588
$ws.Cells.Item(7, 19).Value = 0.001288911522105938
$ws.Cells.Item(7, 20).Value = 0.001276224433172173
$ws.Cells.Item(7, 21).Value = 0.001268880466183995
$ws.Cells.Item(7, 22).Value = 0.001261995234418299
$ws.Cells.Item(7, 23).Value = 0.001256553197283643
$ws.Cells.Item(7, 24).Value = 0.001256553197283643
$ws.Cells.Item(7, 25).Value = 0.001251094378797799

$ws.Cells.Item(8, 3).Value = 0.2917087078094482
$ws.Cells.Item(8, 5).Value = 54.73639758962963
$ws.Cells.Item(8, 6).Value = 0.002010819686020787
$ws.Cells.Item(8, 7).Value = 0.001671325130938397
$ws.Cells.Item(8, 8).Value = 0.001536448147294739
$ws.Cells.Item(8, 9).Value = 0.001498775857768871
$ws.Cells.Item(8, 10).Value = 0.00140620600907169
$ws.Cells.Item(8, 11).Value = 0.001301932025299467
$ws.Cells.Item(8, 12).Value = 0.001289562744469225
$ws.Cells.Item(8, 13).Value = 0.001183495940981106
$ws.Cells.Item(8, 14).Value = 0.00115209442144881
$ws.Cells.Item(8, 15).Value = 0.00115209442144881
$ws.Cells.Item(8, 16).Value = 0.00115209442144881
$ws.Cells.Item(8, 17).Value = 0.001151130349849571
$ws.Cells.Item(8, 18).Value = 0.001151130349849571
$ws.Cells.Item(8, 19).Value = 0.001151130349849571
$ws.Cells.Item(8, 20).Value = 0.001151130349849571
$ws.Cells.Item(8, 21).Value = 0.001140594432344423
$ws.Cells.Item(8, 22).Value = 0.001093654096172814
$ws.Cells.Item(8, 23).Value = 0.001093654096172814
$ws.Cells.Item(8, 24).Value = 0.001084439830985651
$ws.Cells.Item(8, 25).Value = 0.001066986307790051

$ws.Cells.Item(9, 3).Value = 0.2688229084014893
$ws.Cells.Item(9, 5).Value = 54.30725504569273
$ws.Cells.Item(9, 7).Value = 0.001701279647527097
$ws.Cells.Item(9, 8).Value = 0.001669715375980043
$ws.Cells.Item(9, 9).Value = 0.001389192334987158
$ws.Cells.Item(9, 10).Value = 0.001304440730238238
$ws.Cells.Item(9, 11).Value = 0.00126419389355216
$ws.Cells.Item(9, 12).Value = 0.00126419389355216
$ws.Cells.Item(9, 13).Value = 0.00126419389355216
$ws.Cells.Item(9, 14).Value = 0.00119720180547912
$ws.Cells.Item(9, 15).Value = 0.001191414051093862
$ws.Cells.Item(9, 16).Value = 0.001191414051093862
$ws.Cells.Item(9, 17).Value = 0.001188393472400269
$ws.Cells.Item(9, 18).Value = 0.001130031832062439
$ws.Cells.Item(9, 19).Value = 0.001130031832062439
$ws.Cells.Item(9, 20).Value = 0.001117134845701559
$ws.Cells.Item(9, 21).Value = 0.001096553573281507
$ws.Cells.Item(9, 22).Value = 0.001083450178865776
$ws.Cells.Item(9, 23).Value = 0.00106884718162318
$ws.Cells.Item(9, 24).Value = 0.001061182453545326
$ws.Cells.Item(9, 25).Value = 0.001058620956056388

$ws.Cells.Item(10, 3).Value = 0.2356662750244141
$ws.Cells.Item(10, 5).Value = 53.49491020689493
$ws.Cells.Item(10, 6).Value = 0.002130614768142741
$ws.Cells.Item(10, 7).Value = 0.001758391186359664
$ws.Cells.Item(10, 8).Value = 0.001683917744534515
$ws.Cells.Item(10, 9).Value = 0.001563539045679032
$ws.Cells.Item(10, 10).Value = 0.001419065485602971
$ws.Cells.Item(10, 11).Value = 0.001403729042156956
$ws.Cells.Item(10, 12).Value = 0.001386727614262391
$ws.Cells.Item(10, 13).Value = 0.001314304914308257
$ws.Cells.Item(10, 14).Value = 0.001287274446995336
$ws.Cells.Item(10, 15).Value = 0.001204998213667074
$ws.Cells.Item(10, 16).Value = 0.001184539499936576
$ws.Cells.Item(10, 17).Value = 0.001152813295408357
$ws.Cells.Item(10, 18).Value = 0.00108727825076318
$ws.Cells.Item(10, 19).Value = 0.00108727825076318
$ws.Cells.Item(10, 20).Value = 0.001084674675572213
$ws.Cells.Item(10, 21).Value = 0.001063010418777791
$ws.Cells.Item(10, 22).Value = 0.001063010418777791
$ws.Cells.Item(10, 23).Value = 0.001061957676724914
$ws.Cells.Item(10, 24).Value = 0.001056716613344093
$ws.Cells.Item(10, 25).Value = 0.001042785774013546

$ws.Cells.Item(11, 3).Value = 0.2447490692138672
$ws.Cells.Item(11, 5).Value = 58.52000246671923
$ws.Cells.Item(11, 6).Value = 0.002183817042338903
$ws.Cells.Item(11, 7).Value = 0.001793523527633342
$ws.Cells.Item(11, 8).Value = 0.001723264246227391
$ws.Cells.Item(11, 9).Value = 0.001579643025206589
$ws.Cells.Item(11, 10).Value = 0.001551287926628191
$ws.Cells.Item(11, 11).Value = 0.001414899988201469
$ws.Cells.Item(11, 12).Value = 0.001325400678003974
$ws.Cells.Item(11, 13).Value = 0.001325400678003974
$ws.Cells.Item(11, 14).Value = 0.001325400678003974
$ws.Cells.Item(11, 15).Value = 0.001325400678003974
$ws.Cells.Item(11, 16).Value = 0.001298521472862984
$ws.Cells.Item(11, 17).Value = 0.001291218482849366
$ws.Cells.Item(11, 18).Value = 0.001258166678383781
$ws.Cells.Item(11, 19).Value = 0.001189794003234938
$ws.Cells.Item(11, 20).Value = 0.001189794003234938
$ws.Cells.Item(11, 21).Value = 0.001189794003234938
$ws.Cells.Item(11, 22).Value = 0.001189794003234938
$ws.Cells.Item(11, 23).Value = 0.001168898834890059
$ws.Cells.Item(11, 24).Value = 0.001168898834890059
$ws.Cells.Item(11, 25).Value = 0.001140740788824936
